# Applies the per-cell leaderboard refresh described in the commit
# "Code updated 23-05-04 11:02:53": updated Rank/ID/Name/season values
# for the rows whose underlying stats changed between scrapes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = '''56264'
$ws.Range("E2").Value = '''2792'

# Row 3
$ws.Range("A3").Value = '''68294'
$ws.Range("E3").Value = '''2511'

# Row 5
$ws.Range("A5").Value = '''49558'
$ws.Range("E5").Value = '''3130'

# Row 7
$ws.Range("A7").Value = '''44235'
$ws.Range("E7").Value = '''3624'

# Row 8
$ws.Range("A8").Value = '''10269'
$ws.Range("E8").Value = '''5544'

# Row 9
$ws.Range("A9").Value = '''13043'
$ws.Range("E9").Value = '''5336'

# Row 10
$ws.Range("A10").Value = '''17536'
$ws.Range("E10").Value = '''5046'

# Row 11
$ws.Range("A11").Value = '''18857'
$ws.Range("E11").Value = '''4981'

# Row 12
$ws.Range("A12").Value = '''31607'
$ws.Range("E12").Value = '''4390'

# Row 13
$ws.Range("A13").Value = '''52536'

# Row 14
$ws.Range("A14").Value = '''66865'
$ws.Range("E14").Value = '''2530'

# Row 16
$ws.Range("A16").Value = '''13899'
$ws.Range("B16").Value = '''8057001'
$ws.Range("C16").Value = '㊥兵者诡道也'
$ws.Range("E16").Value = '''5275'

# Row 17
$ws.Range("A17").Value = '''16078'
$ws.Range("E17").Value = '''5134'

# Row 18
$ws.Range("A18").Value = '''16303'
$ws.Range("B18").Value = '''31495601'
$ws.Range("C18").Value = '陈晓军'
$ws.Range("E18").Value = '''5119'

# Row 19
$ws.Range("A19").Value = '''16473'
$ws.Range("E19").Value = '''5107'

# Row 20
$ws.Range("A20").Value = '''20644'
$ws.Range("E20").Value = '''4889'

# Row 21
$ws.Range("A21").Value = '''21754'
$ws.Range("E21").Value = '''4831'

# Row 22
$ws.Range("A22").Value = '''21945'
$ws.Range("E22").Value = '''4819'

# Row 23
$ws.Range("A23").Value = '''29718'
$ws.Range("B23").Value = '''3649043'
$ws.Range("C23").Value = 'Dj6106'
$ws.Range("E23").Value = '''4470'

# Row 24
$ws.Range("A24").Value = '''32368'
$ws.Range("B24").Value = '''56585361'
$ws.Range("C24").Value = '"㊥ go策划我要ali"'
$ws.Range("E24").Value = '''4358'

# Row 25
$ws.Range("A25").Value = '''32939'
$ws.Range("B25").Value = '''56732705'
$ws.Range("C25").Value = '时间温柔皆遗憾'
$ws.Range("E25").Value = '''4332'

# Row 26
$ws.Range("A26").Value = '''33224'
$ws.Range("B26").Value = '''58839983'
$ws.Range("C26").Value = '每逢佳节胖六斤'
$ws.Range("E26").Value = '''4320'

# Row 27
$ws.Range("A27").Value = '''37501'
$ws.Range("E27").Value = '''4132'

# Row 28
$ws.Range("A28").Value = '''39296'
$ws.Range("E28").Value = '''4052'

# Row 29
$ws.Range("A29").Value = '''43576'
$ws.Range("E29").Value = '''3701'

# Row 30
$ws.Range("A30").Value = '''6978'
$ws.Range("E30").Value = '''5815'

# Row 31
$ws.Range("A31").Value = '''7739'
$ws.Range("E31").Value = '''5747'

# Row 32
$ws.Range("A32").Value = '''10982'
$ws.Range("E32").Value = '''5495'

# Row 33
$ws.Range("A33").Value = '''12217'
$ws.Range("B33").Value = '''45967307'
$ws.Range("C33").Value = 'Ricky'
$ws.Range("E33").Value = '''5399'

# Row 34
$ws.Range("A34").Value = '''12825'
$ws.Range("B34").Value = '''55317038'
$ws.Range("C34").Value = 'necman12345'
$ws.Range("E34").Value = '''5352'

# Row 35
$ws.Range("A35").Value = '''13201'
$ws.Range("E35").Value = '''5324'

# Row 36
$ws.Range("A36").Value = '''18602'
$ws.Range("E36").Value = '''4993'

# Row 37
$ws.Range("A37").Value = '''20140'
$ws.Range("E37").Value = '''4915'

# Row 38
$ws.Range("A38").Value = '''28352'
$ws.Range("E38").Value = '''4526'

# Row 39
$ws.Range("A39").Value = '''30575'
$ws.Range("E39").Value = '''4435'

# Row 40
$ws.Range("A40").Value = '''31648'
$ws.Range("E40").Value = '''4388'

# Row 41
$ws.Range("A41").Value = '''32163'
$ws.Range("B41").Value = '''52997727'
$ws.Range("C41").Value = 'larios'
$ws.Range("E41").Value = '''4367'

# Row 42
$ws.Range("A42").Value = '''32593'
$ws.Range("B42").Value = '''56379103'
$ws.Range("C42").Value = 'Globalking'
$ws.Range("E42").Value = '''4348'

# Row 43
$ws.Range("A43").Value = '''35021'
$ws.Range("E43").Value = '''4243'

# Row 44
$ws.Range("A44").Value = '''39392'
$ws.Range("B44").Value = '''55634661'
$ws.Range("C44").Value = 'Opalus'
$ws.Range("E44").Value = '''4047'

# Row 45
$ws.Range("A45").Value = '''40447'
$ws.Range("E45").Value = '''3996'

# Row 46
$ws.Range("A46").Value = '''40619'
$ws.Range("B46").Value = '''58203298'
$ws.Range("C46").Value = '权旨qua'
$ws.Range("E46").Value = '''3994'

# Row 47
$ws.Range("A47").Value = '''41848'
$ws.Range("B47").Value = '''32316256'
$ws.Range("C47").Value = '"秋の風 .."'
$ws.Range("E47").Value = '''3909'

# Row 48
$ws.Range("A48").Value = '''47596'
$ws.Range("E48").Value = '''3288'

# Row 49
$ws.Range("A49").Value = '''49000'
$ws.Range("B49").Value = '''48634530'
$ws.Range("C49").Value = 'leezhenrui'
$ws.Range("E49").Value = '''3174'

# Row 50
$ws.Range("A50").Value = '''58240'
$ws.Range("E50").Value = '''2726'

# Row 51
$ws.Range("A51").Value = '''69235'

# Row 52
$ws.Range("A52").Value = '''61902'
$ws.Range("E52").Value = '''2624'

# Row 53
$ws.Range("A53").Value = '''51007'
$ws.Range("E53").Value = '''3030'

# Row 56
$ws.Range("A56").Value = '''41547'
$ws.Range("E56").Value = '''3947'

# Row 57
$ws.Range("A57").Value = '''50965'
$ws.Range("E57").Value = '''3034'

# Row 58
$ws.Range("A58").Value = '''59477'
$ws.Range("E58").Value = '''2687'

# Row 59
$ws.Range("A59").Value = '''69768'
$ws.Range("E59").Value = '''2496'

# Row 60
$ws.Range("A60").Value = '''103635'
$ws.Range("E60").Value = '''1506'

# Row 61
$ws.Range("A61").Value = '''103736'
$ws.Range("E61").Value = '''1504'

# Row 62
$ws.Range("A62").Value = '''110408'
$ws.Range("E62").Value = '''1428'

# Row 63
$ws.Range("A63").Value = '''119310'
$ws.Range("E63").Value = '''1302'

# Row 75
$ws.Range("A75").Value = '''49726'
$ws.Range("E75").Value = '''3118'

# Row 78
$ws.Range("A78").Value = '''89275'
$ws.Range("E78").Value = '''1862'

# Row 79
$ws.Range("A79").Value = '''94993'
$ws.Range("E79").Value = '''1636'

# Row 80
$ws.Range("A80").Value = '''154507'

# Row 81
$ws.Range("A81").Value = '''204086'
